# Update "想去人数" (F column) figures for several exhibition rows on both
# the "展览" and "全部类型" sheets (they carry duplicate data).

$wb = $excel.ActiveWorkbook

$updates = @{
    "F9"  = 6547
    "F11" = 138
    "F12" = 1026
    "F13" = 345
    "F15" = 181
    "F16" = 494
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
